# ============================================================================
# Adds player-info & extended batting-stats scraping output to the
# PlayerPerformance workbook:
#   1. New "Player Info" sheet (inserted first)
#   2. "ODI Batting"  - MATCH_CARD_LINK -> MATCH_CODE (store bare match code,
#      not the full scorecard URL); also drops a stray empty B6 cell.
#   3. "ODI Bowling"  - same MATCH_CARD_LINK -> MATCH_CODE treatment.
#   4. New "ODI Batting Extra" sheet (appended last) with per-match batting
#      detail (batting position, boundaries, % of team runs, MoM award).
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# Helper: stamp a header row with the same look'n'feel the workbook already
# uses for its other header rows (bold, thin box border, centered/top align).
# ----------------------------------------------------------------------------
function Format-HeaderRow($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = -4108   # xlHAlignCenter
    $range.VerticalAlignment = -4160     # xlVAlignTop
    $range.Borders.LineStyle = 1         # xlContinuous
}

# ----------------------------------------------------------------------------
# 1. "Player Info" sheet - inserted before the existing first sheet.
# ----------------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$playerInfo = $wb.Worksheets.Add($odiBatting)
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"
Format-HeaderRow($playerInfo.Range("A1:D1"))

$piIdCell = $playerInfo.Cells.Item(2, 1)
$piIdCell.NumberFormat = "@"
$piIdCell.Value = "4487"
$playerInfo.Range("B2").Value = "Mashtayage Danushka Gunathilaka"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

# ----------------------------------------------------------------------------
# 2. "ODI Batting" - MATCH_CARD_LINK -> MATCH_CODE
# ----------------------------------------------------------------------------
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingLastRow = $battingSheet.Cells.Item(1, 1).End(4).Row
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $link = $cell.Value()
    if ($link -match "MatchCode=(\d+)") {
        $cell.NumberFormat = "@"
        $cell.Value = $matches[1]
    }
}

# Drop the stray empty B6 cell (INNING_NUMBER left blank for the "did not
# bat" match) so the cell is absent rather than an empty placeholder.
$battingSheet.Cells.Item(6, 2).ClearContents()

# ----------------------------------------------------------------------------
# 3. "ODI Bowling" - MATCH_CARD_LINK -> MATCH_CODE
# ----------------------------------------------------------------------------
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingLastRow = $bowlingSheet.Cells.Item(1, 1).End(4).Row
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $link = $cell.Value()
    if ($link -match "MatchCode=(\d+)") {
        $cell.NumberFormat = "@"
        $cell.Value = $matches[1]
    }
}

# ----------------------------------------------------------------------------
# 4. "ODI Batting Extra" - new sheet appended at the end.
# ----------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$extraSheet = $wb.Worksheets.Add($null, $lastSheet)
$extraSheet.Name = "ODI Batting Extra"

$extraSheet.Range("A1").Value = "MATCH_CODE"
$extraSheet.Range("B1").Value = "BATTING_POSITION"
$extraSheet.Range("C1").Value = "NUM_4"
$extraSheet.Range("D1").Value = "NUM_6"
$extraSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extraSheet.Range("F1").Value = "MAN_OF_MATCH"
Format-HeaderRow($extraSheet.Range("A1:F1"))

$extraRows = @(
    @{ Code = "4062"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    Mom = "NO" },
    @{ Code = "4096"; Pos = 1;     N4 = "0";   N6 = "0";   Pct = "0.88%";  Mom = "NO" },
    @{ Code = "4098"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    Mom = "NO" },
    @{ Code = "4099"; Pos = 1;     N4 = "2";   N6 = "0";   Pct = "6.05%";  Mom = "NO" },
    @{ Code = "4122"; Pos = 1;     N4 = "3";   N6 = "2";   Pct = "42.17%"; Mom = "NO" },
    @{ Code = "4124"; Pos = 1;     N4 = "1";   N6 = "0";   Pct = "2.71%";  Mom = "NO" },
    @{ Code = "4231"; Pos = 2;     N4 = "3";   N6 = "0";   Pct = "13.19%"; Mom = "NO" },
    @{ Code = "4232"; Pos = 2;     N4 = "9";   N6 = "0";   Pct = "23.83%"; Mom = "NO" },
    @{ Code = "4233"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    Mom = "NO" },
    @{ Code = "4375"; Pos = 1;     N4 = "2";   N6 = "0";   Pct = "5.88%";  Mom = "NO" },
    @{ Code = "4376"; Pos = 1;     N4 = "16";  N6 = "1";   Pct = "44.78%"; Mom = "NO" },
    @{ Code = "4449"; Pos = 1;     N4 = "7";   N6 = "0";   Pct = "23.71%"; Mom = "NO" },
    @{ Code = "4450"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    Mom = "NO" },
    @{ Code = "4451"; Pos = 1;     N4 = "6";   N6 = "0";   Pct = "13.14%"; Mom = "NO" },
    @{ Code = "4463"; Pos = 1;     N4 = "5";   N6 = "0";   Pct = "9.38%";  Mom = "NO" },
    @{ Code = "4464"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    Mom = "NO" },
    @{ Code = "4465"; Pos = 1;     N4 = "5";   N6 = "1";   Pct = "13.64%"; Mom = "NO" },
    @{ Code = "4594"; Pos = $null; N4 = $null; N6 = $null; Pct = $null;    Mom = "NO" },
    @{ Code = "4597"; Pos = 1;     N4 = "1";   N6 = "0";   Pct = "8.18%";  Mom = "NO" },
    @{ Code = "4603"; Pos = 1;     N4 = "1";   N6 = "0";   Pct = "5.00%";  Mom = "NO" }
)

$r = 2
foreach ($row in $extraRows) {
    $codeCell = $extraSheet.Cells.Item($r, 1)
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $row.Code

    if ($row.Pos -ne $null) {
        # BATTING_POSITION is the one genuinely-numeric column.
        $extraSheet.Cells.Item($r, 2).Value = $row.Pos
    }
    if ($row.N4 -ne $null) {
        $n4Cell = $extraSheet.Cells.Item($r, 3)
        $n4Cell.NumberFormat = "@"
        $n4Cell.Value = $row.N4
    }
    if ($row.N6 -ne $null) {
        $n6Cell = $extraSheet.Cells.Item($r, 4)
        $n6Cell.NumberFormat = "@"
        $n6Cell.Value = $row.N6
    }
    if ($row.Pct -ne $null) {
        $pctCell = $extraSheet.Cells.Item($r, 5)
        $pctCell.NumberFormat = "@"
        $pctCell.Value = $row.Pct
    }

    $momCell = $extraSheet.Cells.Item($r, 6)
    $momCell.NumberFormat = "@"
    $momCell.Value = $row.Mom

    $r = $r + 1
}

# Keep the first sheet ("Player Info") as the active tab, as in the source
# workbook (activeTab 0).
$playerInfo.Activate()
